# Trade #6 closed at 2026-02-17 12:27:04 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.01
$summary.Range("B4").Value = 0.01
$summary.Range("B5").Value = 0.03
$summary.Range("B6").Value = 6
$summary.Range("B7").Value = 3
$summary.Range("B9").Value = 50

# --- Strategy Status sheet (MarketMaking strategy row, row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.01
$status.Range("D4").Value = 6
$status.Range("E4").Value = 0.01
$status.Range("F4").Value = 0.01
$status.Range("G4").Value = 50

# --- New trade row (trade #6) appended to "All Trades" and "MarketMaking" sheets ---
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 7

    $ws.Cells.Item($row, 1).Value = 6

    # "2026-02-17" must stay plain text (like the rows above it) instead of
    # being auto-converted to a date serial number by Excel's type inference.
    # Force the cell to Text format, write the string, then strip the
    # number-format back off so the cell ends up identical in shape to its
    # neighbours (plain text, default style).
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).Value = "12:26:58"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.86
    $ws.Cells.Item($row, 7).Value = 0.88
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 2.3256
    $ws.Cells.Item($row, 10).Value = 0.02
    $ws.Cells.Item($row, 11).Value = 100.01
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}
